# TC03 workbook: add a leading "TabName" column (with a "CasesTab" tag row)
# and refresh the Neo4j queries used to build the trial list / stat sheets
# for the new gender + ethnicity CTDC test case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing four columns (query/StatQuery/dbExcel/WebExcel file names)
# one column to the right so a new "TabName" column can be inserted at A.
$ws.Columns.Item(1).Insert()

# New column A: tab-name metadata used by the automation runner.
$ws.Range("A1").Value2 = "TabName"
$ws.Range("A2").Value2 = "CasesTab"

# Updated Cypher queries (column B = main trial/case query, column C = stat query).
$ws.Range("B2").Value2 = 'MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)
    WHERE c.ethnicity IN [''UNKNOWN'']
WITH DISTINCT c, a, ct
RETURN 
    COALESCE(c.case_id, '''') AS `Case ID`,
    COALESCE(ct.clinical_trial_designation, '''') AS `Trial Code`,
    COALESCE(a.arm_id, '''') AS `Arm`,
    COALESCE(a.arm_drug, '''') AS `Arm Treatment`,
    COALESCE(c.disease, '''') AS `Diagnosis`,
    COALESCE(c.gender, '''') AS `Gender`,
    COALESCE(c.race, '''') AS `Race`,
    COALESCE(c.ethnicity, '''') AS `Ethnicity`'
$ws.Range("C2").Value2 = 'MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)
    WHERE WHERE c.ethnicity IN [''UNKNOWN'']
OPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)
RETURN 
    COUNT(DISTINCT f) AS number_of_files,
    COUNT(DISTINCT c.case_id) AS number_of_cases,
    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials'

# Column A is a short bestfit-style column; keep the rest of the columns as-is
# (they were simply shifted right by the insert above).
$ws.Columns.Item(1).ColumnWidth = 8

# Row 2 grows taller because the wrapped query text is now longer.
$ws.Rows.Item(2).RowHeight = 174

# Selection moves onto the (now relocated) trial query cell.
$ws.Range("B2").Select() | Out-Null
